$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (before the existing ST_Load column), shifting
# ST_Load, ST_Temp, RT_Load, RT_Temp, WeightOnBit, DrillingVelocity one
# column to the right.
$ws.Columns("H:H").Insert()

# Populate the header of the newly inserted column.
$ws.Cells.Item(1, 8).Value = "DrillBit_Rotation"

# Match the author's final selection (the newly added cell).
$ws.Range("H1").Select() | Out-Null
